# [Kadastro App] Yeni kayit eklendi: 2955
#
# Appends one new record row (row 46) to both the master "Kayitlar" sheet
# and the per-district "Erdemli" sheet, which mirrors the same records.
# Every column is stored as text (matching the rest of the table, which
# uses t="str" cells even for numeric-looking values like the record id,
# date and parcel count), so numeric/date-looking values are entered with
# a leading apostrophe to force text, then the quote-prefix visual style
# is cleared so the cell keeps its plain "Normal" look.

$wb = $excel.ActiveWorkbook

$newRecord = @{
    A = "2955"
    B = "2025-09-09"
    C = "Erdemli"
    D = "1"
    E = "ÇAP"
    F = "CEMAL TİMUROĞLU (K.Teknisyeni)"
}

$targetSheets = @("Kayitlar", "Erdemli")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newRow = $ws.UsedRange.Rows.Count() + 1

    foreach ($col in @("A", "B", "C", "D", "E", "F")) {
        $cell = $ws.Range($col + $newRow)
        $cell.Value = "'" + $newRecord[$col]
        $cell.Style = "Normal"
    }
}
